$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("production_tax_credits")
$ws2 = $wb.Worksheets.Item("investment_tax_credits")

# --- production_tax_credits: update discount-rate inputs (E2:E9) with the
# actual computed discount rates; downstream formulas in F/G/I recalc
# automatically. Rows 10-18 keep their original input values.
$ws1.Range("E2").Value = 0.0320752195121951
$ws1.Range("E3").Value = 0.0564730561021376
$ws1.Range("E4").Value = 0.051519516365778
$ws1.Range("E5").Value = 0.0538314857296738
$ws1.Range("E6").Value = 0.0438437157985803
$ws1.Range("E7").Value = 0.0519007613262936
$ws1.Range("E8").Value = 0.0515227657596506
$ws1.Range("E9").Value = 0.0515227657596506

# New discount rates get a dedicated 3-decimal display format.
$ws1.Range("E2:E9").NumberFormat = "0.000"

# --- Selections / active sheet bookkeeping -------------------------------
[void]$ws1.Range("I2").Select()
[void]$ws2.Range("G9").Select()

# investment_tax_credits becomes the active tab.
$ws2.Activate()
